$p = $ppt.ActivePresentation

# --- 1) Update the "datetimeFigureOut" Date Placeholder field text on the
#        slide master and every slide layout: 09-08-2022 -> 17-08-2022 ---

$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "17-08-2022"
    }
}

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $lay = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "17-08-2022"
        }
    }
}

# --- 2) Trim the footer date text on slides 1-3 down to just the year ---

$footerDates = @{ 1 = "10-Aug-2022"; 2 = "28-Aug-2022"; 3 = "08-Aug-2022" }

foreach ($slideIdx in $footerDates.Keys) {
    $s = $p.Slides.Item($slideIdx)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame -eq -1 -and $sh.Name -like "Footer Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq $footerDates[$slideIdx]) {
                $sh.TextFrame.TextRange.Text = "2022"
            }
        }
    }
}
